$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E width (new column 5); closest achievable to 10.7109375 via the COM emulation is integer ColumnWidth=10 ---
$ws.Columns.Item(5).ColumnWidth = 10

# --- Row 1 / Row 2: new E1/E2 empty-text cells matching style of B1:D1 / B2:D2 (s=4: right-aligned) ---
# Using a formula-then-freeze trick so the empty string is stored as an explicit text value
# (not a fully blank cell), matching the existing B1/C1/D1 encoding.
$ws.Range("E1").Formula = "="""""
$ws.Range("E1").Value = $ws.Range("E1").Value
$ws.Range("E1").HorizontalAlignment = -4152
$ws.Range("E2").Formula = "="""""
$ws.Range("E2").Value = $ws.Range("E2").Value
$ws.Range("E2").HorizontalAlignment = -4152

# --- Row 3: existing B3:E3 date cells get new date strings (safe - not brand new cells) ---
$ws.Range("B3").Value = "28-08-2024"
$ws.Range("C3").Value = "29-08-2024"
$ws.Range("D3").Value = "30-08-2024"
$ws.Range("E3").Value = "31-08-2024"

# --- Row 3: new F3 date cell. A direct .Value assignment of "01-09-2024" on a brand-new cell
# gets auto-recognized as a real date (ambiguous D-M / M-D) and stored as a date serial.
# Route the literal text through a formula first, then freeze it back to a value, which keeps it text;
# then reapply the s=3 look (right align + bold) in the same order used elsewhere to avoid a stray style.
$ws.Range("F3").Formula = "=""01-09-2024"""
$ws.Range("F3").Value = $ws.Range("F3").Value
$ws.Range("F3").HorizontalAlignment = -4152
$ws.Range("F3").Font.Bold = $true

# --- Rows 4-27: existing B:E cells simply get new values ---
$ws.Range("B4").Value = "410,00"
$ws.Range("C4").Value = "413,70"
$ws.Range("D4").Value = "522,00"
$ws.Range("E4").Value = "482,00"
$ws.Range("B5").Value = "388,00"
$ws.Range("C5").Value = "384,73"
$ws.Range("D5").Value = "450,96"
$ws.Range("E5").Value = "450,00"
$ws.Range("B6").Value = "376,80"
$ws.Range("C6").Value = "367,53"
$ws.Range("D6").Value = "427,60"
$ws.Range("E6").Value = "429,60"
$ws.Range("B7").Value = "370,80"
$ws.Range("C7").Value = "376,71"
$ws.Range("D7").Value = "425,00"
$ws.Range("E7").Value = "418,90"
$ws.Range("B8").Value = "383,30"
$ws.Range("C8").Value = "390,00"
$ws.Range("D8").Value = "440,00"
$ws.Range("E8").Value = "413,23"
$ws.Range("B9").Value = "402,00"
$ws.Range("C9").Value = "415,96"
$ws.Range("D9").Value = "490,00"
$ws.Range("E9").Value = "457,80"
$ws.Range("B10").Value = "520,00"
$ws.Range("C10").Value = "550,00"
$ws.Range("D10").Value = "657,31"
$ws.Range("E10").Value = "481,43"
$ws.Range("B11").Value = "530,00"
$ws.Range("C11").Value = "572,36"
$ws.Range("D11").Value = "700,00"
$ws.Range("E11").Value = "492,00"
$ws.Range("B12").Value = "479,10"
$ws.Range("C12").Value = "478,89"
$ws.Range("D12").Value = "570,00"
$ws.Range("E12").Value = "460,00"
$ws.Range("B13").Value = "382,60"
$ws.Range("C13").Value = "378,50"
$ws.Range("D13").Value = "455,79"
$ws.Range("E13").Value = "353,20"
$ws.Range("B14").Value = "259,99"
$ws.Range("C14").Value = "230,00"
$ws.Range("D14").Value = "349,23"
$ws.Range("E14").Value = "190,00"
$ws.Range("B15").Value = "150,00"
$ws.Range("C15").Value = "150,00"
$ws.Range("D15").Value = "310,11"
$ws.Range("E15").Value = "65,30"
$ws.Range("B16").Value = "70,00"
$ws.Range("C16").Value = "85,99"
$ws.Range("D16").Value = "273,99"
$ws.Range("E16").Value = "9,36"
$ws.Range("B17").Value = "55,00"
$ws.Range("C17").Value = "50,00"
$ws.Range("D17").Value = "277,99"
$ws.Range("E17").Value = "0,00"
$ws.Range("B18").Value = "90,00"
$ws.Range("C18").Value = "108,88"
$ws.Range("D18").Value = "319,20"
$ws.Range("E18").Value = "1,00"
$ws.Range("B19").Value = "282,89"
$ws.Range("C19").Value = "288,00"
$ws.Range("D19").Value = "398,00"
$ws.Range("E19").Value = "211,99"
$ws.Range("B20").Value = "388,49"
$ws.Range("C20").Value = "400,00"
$ws.Range("D20").Value = "552,99"
$ws.Range("E20").Value = "411,00"
$ws.Range("B21").Value = "508,93"
$ws.Range("C21").Value = "627,30"
$ws.Range("D21").Value = "750,00"
$ws.Range("E21").Value = "516,60"
$ws.Range("B22").Value = "769,99"
$ws.Range("C22").Value = "1050,00"
$ws.Range("D22").Value = "1471,00"
$ws.Range("E22").Value = "778,40"
$ws.Range("B23").Value = "1050,00"
$ws.Range("C23").Value = "1600,00"
$ws.Range("D23").Value = "2071,18"
$ws.Range("E23").Value = "995,96"
$ws.Range("B24").Value = "1050,00"
$ws.Range("C24").Value = "1511,53"
$ws.Range("D24").Value = "2000,00"
$ws.Range("E24").Value = "867,40"
$ws.Range("B25").Value = "540,00"
$ws.Range("C25").Value = "749,77"
$ws.Range("D25").Value = "824,89"
$ws.Range("E25").Value = "668,50"
$ws.Range("B26").Value = "461,10"
$ws.Range("C26").Value = "590,00"
$ws.Range("D26").Value = "561,70"
$ws.Range("E26").Value = "531,83"
$ws.Range("B27").Value = "400,00"
$ws.Range("C27").Value = "488,60"
$ws.Range("D27").Value = "480,00"
$ws.Range("E27").Value = "466,55"

# --- Rows 4-27: new F column cells; set value then right-align to match style s=4 ---
$ws.Range("F4").Value = "450,00"
$ws.Range("F4").HorizontalAlignment = -4152
$ws.Range("F5").Value = "420,00"
$ws.Range("F5").HorizontalAlignment = -4152
$ws.Range("F6").Value = "396,00"
$ws.Range("F6").HorizontalAlignment = -4152
$ws.Range("F7").Value = "392,00"
$ws.Range("F7").HorizontalAlignment = -4152
$ws.Range("F8").Value = "398,00"
$ws.Range("F8").HorizontalAlignment = -4152
$ws.Range("F9").Value = "403,70"
$ws.Range("F9").HorizontalAlignment = -4152
$ws.Range("F10").Value = "400,00"
$ws.Range("F10").HorizontalAlignment = -4152
$ws.Range("F11").Value = "397,60"
$ws.Range("F11").HorizontalAlignment = -4152
$ws.Range("F12").Value = "306,90"
$ws.Range("F12").HorizontalAlignment = -4152
$ws.Range("F13").Value = "80,00"
$ws.Range("F13").HorizontalAlignment = -4152
$ws.Range("F14").Value = "0,01"
$ws.Range("F14").HorizontalAlignment = -4152
$ws.Range("F15").Value = "-53,50"
$ws.Range("F15").HorizontalAlignment = -4152
$ws.Range("F16").Value = "-96,56"
$ws.Range("F16").HorizontalAlignment = -4152
$ws.Range("F17").Value = "-175,49"
$ws.Range("F17").HorizontalAlignment = -4152
$ws.Range("F18").Value = "-131,37"
$ws.Range("F18").HorizontalAlignment = -4152
$ws.Range("F19").Value = "-0,02"
$ws.Range("F19").HorizontalAlignment = -4152
$ws.Range("F20").Value = "297,99"
$ws.Range("F20").HorizontalAlignment = -4152
$ws.Range("F21").Value = "465,40"
$ws.Range("F21").HorizontalAlignment = -4152
$ws.Range("F22").Value = "597,00"
$ws.Range("F22").HorizontalAlignment = -4152
$ws.Range("F23").Value = "700,00"
$ws.Range("F23").HorizontalAlignment = -4152
$ws.Range("F24").Value = "700,00"
$ws.Range("F24").HorizontalAlignment = -4152
$ws.Range("F25").Value = "542,00"
$ws.Range("F25").HorizontalAlignment = -4152
$ws.Range("F26").Value = "499,98"
$ws.Range("F26").HorizontalAlignment = -4152
$ws.Range("F27").Value = "444,23"
$ws.Range("F27").HorizontalAlignment = -4152
